# Insert a new "Match ID" column at the very start of the sheet (column A),
# pushing all existing columns one position to the right (A->B, ..., AC->AD).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

# Header label for the new column (row 3 is the visible header row).
$ws.Range("A3").Value = "Match ID"

# Fill the new column with the match id (29) for every visible data row.
for ($r = 4; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = 29
}

# Row 20 is a hidden "totals" row; temporarily unhide it so writing the
# value doesn't disturb its row formatting, then re-hide it.
$ws.Rows.Item(20).Hidden = $false
$ws.Cells.Item(20, 1).Value = 29
$ws.Rows.Item(20).Hidden = $true

# Bold the new column's header + data cells (matches the new cell style
# used for "Match ID" / its values).
$ws.Range("A3:A19").Font.Bold = $true

# Leave the new column selected, as in the saved workbook.
$ws.Range("A3:A19").Select() | Out-Null
